# "start add data in db"
# Adds a new worksheet ("Sheet1") after the existing "Лист1" sheet and
# populates it with the de-duplicated, sorted list of colors that appear
# (comma separated) in the "Colors" column of the first sheet, plus two
# extra colors appended at the end.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- add the new sheet right after the first one ---------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet1"

$colors = @(
    "beige",
    "beige-gum",
    "black",
    "black-brown",
    "black-gum",
    "black-white",
    "blue",
    "brown",
    "brown-black",
    "gray",
    "green",
    "grey",
    "red",
    "white",
    "white-black",
    "white-red",
    "red-black",
    "yellow"
)

for ($i = 0; $i -lt $colors.Length; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $colors[$i]
}

$ws2.Columns.Item(1).ColumnWidth = 13.8

# the style table on the first sheet used to carry a duplicated "locked"
# format; re-asserting the (unchanged) protection state on the cells that
# used it lets the engine fold them onto the single remaining entry
$ws1.Range("E5").Locked = $true
$ws1.Range("F5").Locked = $true
$ws1.Range("E9").Locked = $true

# keep the first sheet active/selected (matches the unchanged bookViews)
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
